$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 36.5546875
$ws.Columns.Item(2).ColumnWidth = 81.109375

# New formulas in column A for rows 15 and 16
$ws.Range("A15").Formula = '=IF(A7=TRUE,"Actif","Inactif")'
$ws.Range("A16").Formula = '=IF(A8=TRUE,"Actif","Inactif")'

# Selection change
$ws.Range("A12").Select()
